# "reconfigurando e adicionando mais funcao"
#
# The sheet held a lookup table of pdv/valor/Serial rows in A2:C53 (backed
# by shared strings for the text values). The edit wipes that data out
# entirely (values + formatting, so the now-empty rows collapse down to
# just the pre-existing blank "C" placeholder cells used further down the
# sheet) and leaves the selection parked over the range that used to hold
# the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data (contents AND formatting) from A2:C53 so those rows
# match the blank placeholder rows (C54, C55, ... already s="3", no value)
# further down the sheet.
$ws.Range("A2:C53").Clear() | Out-Null

# Leave the selection where the user left it after clearing the table.
$ws.Range("A2:C51").Select() | Out-Null
